$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '27.654.77'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '1.588.15'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '207.45'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '22.23'
$ws.Range('E8').Value = '  -4.45%  '
$ws.Range('D9').Value = '0.253'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '0.0866'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').Value = '1.813.43'
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').Value = '1.587.13'
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').Value = '0.531'
$ws.Range('E15').Value = '  -4.56%  '
$ws.Range('D16').Value = '27.650.39'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '63.40'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '219.64'
$ws.Range('E18').Value = '  -3.66%  '
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('E22').Value = '  -4.69%  '
$ws.Range('D23').Value = '9.68'
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  -3.46%  '
$ws.Range('D26').Value = '6.82'
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('D28').Value = '15.12'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('E29').Value = '  -4.76%  '
$ws.Range('E30').Value = '  -2.28%  '
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('E32').Value = '  -5.14%  '
$ws.Range('D33').Value = '1.374.13'
$ws.Range('E33').Value = '  -2.99%  '
$ws.Range('E34').Value = '  -5.39%  '
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('D36').Value = '0.974'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('D39').Value = '0.535'
$ws.Range('E39').Value = '  -3.05%  '
$ws.Range('E40').Value = '  -3.13%  '
$ws.Range('E42').Value = '  -3.40%  '
$ws.Range('D43').Value = '64.20'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').Value = '1.724.38'
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('E47').Value = '  -5.15%  '
$ws.Range('D48').Value = '87.56'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('D50').Value = '0.0967'
$ws.Range('E50').Value = '  -4.08%  '
$ws.Range('E51').Value = '  -1.53%  '
